# Update the Alvearie FHIR StructureDefinition workbook to the
# LinuxForHealth fork: new URL, version, publication date and publisher
# on the "Metadata" sheet, plus two follow-on fixes on the "Elements"
# sheet (the ele-1/ext-1 constraint belongs on Extension.extension only,
# and the Fixed Value for Extension.url reflects the new canonical URL).

$wb = $excel.ActiveWorkbook

$newUrl = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-birth-date"

# --- "Metadata" sheet (Property / Value pairs) ---
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = $newUrl
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" sheet (StructureDefinition element table) ---
$elements = $wb.Worksheets.Item("Elements")

# Row 2 = "Extension" (the root element): the ele-1/ext-1 constraint
# previously duplicated here should only live on Extension.extension
# (row 4), so clear it from the root row's Constraint(s) column.
$elements.Range("AI2").Value = ""

# Row 5 = "Extension.url": its Fixed Value now points at the new
# canonical URL for this extension.
$elements.Range("Q5").Value = $newUrl
